$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new value looks like a plain number need to be
# forced to Text format first, otherwise Excel auto-converts the literal
# "26.21" style price string into a numeric value (losing the original
# text representation), same as typing it into a Text-less cell in the UI.
$ws.Range("D2").Value = "26.223.28"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").Value = "1.604.90"
$ws.Range("E3").Value = "  -0.44%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.16"
$ws.Range("E5").Value = "  -0.41%  "
$ws.Range("E6").Value = "  -0.16%  "
$ws.Range("E7").Value = "  +0.59%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0614"
$ws.Range("E9").Value = "  -0.85%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.15"
$ws.Range("E10").Value = "  -1.32%  "
$ws.Range("E11").Value = "  +1.36%  "
$ws.Range("D12").Value = "1.826.49"
$ws.Range("E12").Value = "  -0.57%  "
$ws.Range("D13").Value = "1.598.98"
$ws.Range("E13").Value = "  -0.72%  "
$ws.Range("E14").Value = "  -0.92%  "
$ws.Range("E15").Value = "  +1.00%  "
$ws.Range("D16").Value = "26.214.84"
$ws.Range("E16").Value = "  +0.24%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.33"
$ws.Range("E17").Value = "  +0.77%  "
$ws.Range("E18").Value = "  -0.25%  "
$ws.Range("E19").Value = "  -0.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "203.61"
$ws.Range("E20").Value = "  +2.56%  "
$ws.Range("E21").Value = "  +0.37%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.27"
$ws.Range("E22").Value = "  -2.24%  "
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.94"
$ws.Range("E24").Value = "  +11.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.68"
$ws.Range("E25").Value = "  +1.46%  "
$ws.Range("E27").Value = "  -6.53%  "
$ws.Range("E28").Value = "  +0.11%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.55"
$ws.Range("E29").Value = "  +0.54%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0492"
$ws.Range("E30").Value = "  +3.48%  "
$ws.Range("E31").Value = "  -0.62%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.16"
$ws.Range("E32").Value = "  +0.25%  "
$ws.Range("E33").Value = "  -3.93%  "
$ws.Range("E34").Value = "  -2.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.35"
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").Value = "1.138.11"
$ws.Range("E36").Value = "  +2.69%  "
$ws.Range("E37").Value = "  +5.78%  "
$ws.Range("E38").Value = "  -0.11%  "
$ws.Range("E39").Value = "  -0.52%  "
$ws.Range("E40").Value = "  -0.88%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.495"
$ws.Range("E41").Value = "  -2.27%  "
$ws.Range("E42").Value = "  -1.58%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.21"
$ws.Range("E43").Value = "  +0.72%  "
$ws.Range("D44").Value = "1.740.34"
$ws.Range("E44").Value = "  -0.52%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.03"
$ws.Range("E45").Value = "  -1.25%  "
$ws.Range("E46").Value = "  -2.92%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "54.21"
$ws.Range("E47").Value = "  +0.54%  "
$ws.Range("E48").Value = "  -0.71%  "
$ws.Range("E49").Value = "  -0.18%  "
$ws.Range("E50").Value = "  -0.05%  "
$ws.Range("D51").Value = "0.0₇0953"
$ws.Range("E51").Value = "  -11.84%  "
